# Insert a new row at position 5, shifting existing rows 5..85 down to 6..86,
# and populate the newly inserted row 5 with the new data record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C5").Value = "Arica y Parinacota"
$ws.Range("D5").Value = 44616
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100109
$ws.Range("H5").Value = "Uva"
$ws.Range("I5").Value = 100109001
$ws.Range("J5").Value = "Uva"
$ws.Range("K5").Value = "Thompson seedless"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 270
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 14000
$ws.Range("P5").Value = 13500
$ws.Range("Q5").Value = "`$/caja 18 kilos"
$ws.Range("R5").Value = "Región de Coquimbo"
$ws.Range("S5").Value = 750
$ws.Range("T5").Value = 18
